$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shape = $d.Shapes.Item($i)
    $tr = $shape.TextFrame.TextRange
    $find = $tr.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute("стоп", $false, $false, $false, $false, $false, $true, 1, $false, "конец", 2)
}
